# Add "Start Date" (G) and "Expire Date" (H) columns to the content sheet,
# matching the header styling already used for the other header cells,
# size the new columns, add explanatory header comments (mirroring the
# existing B1/C1 comments) and leave the selection on the new last header
# cell (H1) as the author did.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells -------------------------------------------------
$ws.Range("G1").Value = "Start Date"
$ws.Range("H1").Value = "Expire Date"

# Copy the formatting of the last existing header cell (F1, "Check Cell"
# style) onto the two new header cells so they render identically to the
# rest of the header row.
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column sizing ------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 20.5
$ws.Columns.Item(8).ColumnWidth = 27.83
$ws.Columns.Item(9).ColumnWidth = 11.17

# --- Header cell comments (author: mohamed mahmoud) ---------------------
$ws.Range("G1").AddComment("mohamed mahmoud:`n2020-10-22")
$ws.Range("H1").AddComment("mohamed mahmoud:`n2020-10-27")

# --- Selection / view ----------------------------------------------------
$null = $ws.Range("H1").Select()
